$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 32.935331
$ws.Range("H2").Value = 98.805993
$ws.Range("I2").Value = 0.1836164637112342
$ws.Range("J2").Value = 0.1836164637112342
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.131233999999999
$ws.Range("N2").Value = 24.393702
$ws.Range("O2").Value = 0.02090995573015822
$ws.Range("P2").Value = 0.02090995573015823
$ws.Range("Q2").Value = 267.8048832284539
$ws.Range("R2").Value = 2410.243949056086
$ws.Range("S2").Value = 0.00383941212753011
$ws.Range("T2").Value = 0.003839412127530111

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 32.935331
$ws.Range("H3").Value = 98.805993
$ws.Range("I3").Value = 0.1836164637112342
$ws.Range("J3").Value = 0.1836164637112342
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 243.3763986666667
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.625857000534647
$ws.Range("P3").Value = 0.6258570005346471
$ws.Range("Q3").Value = 8015.682247674626
$ws.Range("R3").Value = 72141.14022907164
$ws.Range("S3").Value = 0.1149176492270919
$ws.Range("T3").Value = 0.1149176492270919

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 32.935331
$ws.Range("H4").Value = 98.805993
$ws.Range("I4").Value = 0.1836164637112342
$ws.Range("J4").Value = 0.1836164637112342
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.9426383333333
$ws.Range("N4").Value = 311.827915
$ws.Range("O4").Value = 0.2672947262403034
$ws.Range("P4").Value = 0.2672947262403035
$ws.Range("Q4").Value = 3423.385198521622
$ws.Range("R4").Value = 30810.4667866946
$ws.Range("S4").Value = 0.04907971240090694
$ws.Range("T4").Value = 0.04907971240090695

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 32.935331
$ws.Range("H5").Value = 98.805993
$ws.Range("I5").Value = 0.1836164637112342
$ws.Range("J5").Value = 0.1836164637112342
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.41874933333333
$ws.Range("N5").Value = 100.256248
$ws.Range("O5").Value = 0.08593831749489127
$ws.Range("P5").Value = 0.08593831749489128
$ws.Range("Q5").Value = 1100.657570899363
$ws.Range("R5").Value = 9905.918138094265
$ws.Range("S5").Value = 0.01577968995570522
$ws.Range("T5").Value = 0.01577968995570523

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 95.562134
$ws.Range("H6").Value = 286.686402
$ws.Range("I6").Value = 0.5327646808765668
$ws.Range("J6").Value = 0.5327646808765667
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.131233999999999
$ws.Range("N6").Value = 24.393702
$ws.Range("O6").Value = 0.02090995573015822
$ws.Range("P6").Value = 0.02090995573015823
$ws.Range("Q6").Value = 777.0380730933559
$ws.Range("R6").Value = 6993.342657840203
$ws.Range("S6").Value = 0.01114008589172089
$ws.Range("T6").Value = 0.01114008589172089

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 95.562134
$ws.Range("H7").Value = 286.686402
$ws.Range("I7").Value = 0.5327646808765668
$ws.Range("J7").Value = 0.5327646808765667
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 243.3763986666667
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.625857000534647
$ws.Range("P7").Value = 0.6258570005346471
$ws.Range("Q7").Value = 23257.56802182142
$ws.Range("R7").Value = 209318.1121963928
$ws.Range("S7").Value = 0.3334345051642065
$ws.Range("T7").Value = 0.3334345051642065

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 95.562134
$ws.Range("H8").Value = 286.686402
$ws.Range("I8").Value = 0.5327646808765668
$ws.Range("J8").Value = 0.5327646808765667
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.9426383333333
$ws.Range("N8").Value = 311.827915
$ws.Range("O8").Value = 0.2672947262403034
$ws.Range("P8").Value = 0.2672947262403035
$ws.Range("Q8").Value = 9932.980332723537
$ws.Range("R8").Value = 89396.82299451184
$ws.Range("S8").Value = 0.1424051895254045
$ws.Range("T8").Value = 0.1424051895254045

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 95.562134
$ws.Range("H9").Value = 286.686402
$ws.Range("I9").Value = 0.5327646808765668
$ws.Range("J9").Value = 0.5327646808765667
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.41874933333333
$ws.Range("N9").Value = 100.256248
$ws.Range("O9").Value = 0.08593831749489127
$ws.Range("P9").Value = 0.08593831749489128
$ws.Range("Q9").Value = 3193.567001904411
$ws.Range("R9").Value = 28742.1030171397
$ws.Range("S9").Value = 0.04578490029523483
$ws.Range("T9").Value = 0.04578490029523482

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.06916133333333
$ws.Range("H10").Value = 111.207484
$ws.Range("I10").Value = 0.2066628180165514
$ws.Range("J10").Value = 0.2066628180165514
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.131233999999999
$ws.Range("N10").Value = 24.393702
$ws.Range("O10").Value = 0.02090995573015822
$ws.Range("P10").Value = 0.02090995573015823
$ws.Range("Q10").Value = 301.4180249850853
$ws.Range("R10").Value = 2712.762224865768
$ws.Range("S10").Value = 0.004321310375795836
$ws.Range("T10").Value = 0.004321310375795836

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 37.06916133333333
$ws.Range("H11").Value = 111.207484
$ws.Range("I11").Value = 0.2066628180165514
$ws.Range("J11").Value = 0.2066628180165514
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 243.3763986666667
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.625857000534647
$ws.Range("P11").Value = 0.6258570005346471
$ws.Range("Q11").Value = 9021.758986900319
$ws.Range("R11").Value = 81195.83088210288
$ws.Range("S11").Value = 0.1293413714058765
$ws.Range("T11").Value = 0.1293413714058765

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 37.06916133333333
$ws.Range("H12").Value = 111.207484
$ws.Range("I12").Value = 0.2066628180165514
$ws.Range("J12").Value = 0.2066628180165514
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.9426383333333
$ws.Range("N12").Value = 311.827915
$ws.Range("O12").Value = 0.2672947262403034
$ws.Range("P12").Value = 0.2672947262403035
$ws.Range("Q12").Value = 3853.066429790651
$ws.Range("R12").Value = 34677.59786811586
$ws.Range("S12").Value = 0.05523988136578376
$ws.Range("T12").Value = 0.05523988136578376

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 37.06916133333333
$ws.Range("H13").Value = 111.207484
$ws.Range("I13").Value = 0.2066628180165514
$ws.Range("J13").Value = 0.2066628180165514
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.41874933333333
$ws.Range("N13").Value = 100.256248
$ws.Range("O13").Value = 0.08593831749489127
$ws.Range("P13").Value = 0.08593831749489128
$ws.Range("Q13").Value = 1238.805010595559
$ws.Range("R13").Value = 11149.24509536003
$ws.Range("S13").Value = 0.01776025486909533
$ws.Range("T13").Value = 0.01776025486909533

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.80362366666667
$ws.Range("H14").Value = 41.410871
$ws.Range("I14").Value = 0.07695603739564764
$ws.Range("J14").Value = 0.07695603739564763
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.131233999999999
$ws.Range("N14").Value = 24.393702
$ws.Range("O14").Value = 0.02090995573015822
$ws.Range("P14").Value = 0.02090995573015823
$ws.Range("Q14").Value = 112.2404940816047
$ws.Range("R14").Value = 1010.164446734442
$ws.Range("S14").Value = 0.001609147335111393
$ws.Range("T14").Value = 0.001609147335111393

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.80362366666667
$ws.Range("H15").Value = 41.410871
$ws.Range("I15").Value = 0.07695603739564764
$ws.Range("J15").Value = 0.07695603739564763
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 243.3763986666667
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.625857000534647
$ws.Range("P15").Value = 0.6258570005346471
$ws.Range("Q15").Value = 3359.476216543302
$ws.Range("R15").Value = 30235.28594888972
$ws.Range("S15").Value = 0.04816347473747216
$ws.Range("T15").Value = 0.04816347473747216

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.80362366666667
$ws.Range("H16").Value = 41.410871
$ws.Range("I16").Value = 0.07695603739564764
$ws.Range("J16").Value = 0.07695603739564763
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.9426383333333
$ws.Range("N16").Value = 311.827915
$ws.Range("O16").Value = 0.2672947262403034
$ws.Range("P16").Value = 0.2672947262403035
$ws.Range("Q16").Value = 1434.785062473774
$ws.Range("R16").Value = 12913.06556226397
$ws.Range("S16").Value = 0.02056994294820819
$ws.Range("T16").Value = 0.02056994294820819

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.80362366666667
$ws.Range("H17").Value = 41.410871
$ws.Range("I17").Value = 0.07695603739564764
$ws.Range("J17").Value = 0.07695603739564763
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.41874933333333
$ws.Range("N17").Value = 100.256248
$ws.Range("O17").Value = 0.08593831749489127
$ws.Range("P17").Value = 0.08593831749489128
$ws.Range("Q17").Value = 461.2998392080008
$ws.Range("R17").Value = 4151.698552872008
$ws.Range("S17").Value = 0.006613472374855892
$ws.Range("T17").Value = 0.006613472374855892
